# "Generate Report for Handback" - refresh the localization-status report
# after a successful handback: flip status text, stamp the new handback
# timestamps, clear the stale "out of date" error details, and widen a
# couple of columns so the (now longer/shorter) text fits.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: widen the zh-cn / de-de status columns (E, F) and
# flip their status text the same way as the per-language sheets ---
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Columns.Item(5).ColumnWidth = 29.1666666666667
$ovw.Columns.Item(6).ColumnWidth = 29.1666666666667

$ovw.Range("E2").Value = "Handed back: in sync with en-US"
$ovw.Range("F2").Value = "Handed back: in sync with en-US"
$ovw.Range("E3").Value = "Handed back: in sync with en-US"
$ovw.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")

# Widen "Status" (C) now that the text is longer, shrink "Error Detail" (P)
# now that it is cleared out.
$zh.Columns.Item(3).ColumnWidth = 29.1666666666667
$zh.Columns.Item(16).ColumnWidth = 12.8333333333333

# Row 2 (a.md) and row 3 (b.md): handback completed, in sync with en-US.
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "Handed back: in sync with en-US"

# Fresh handback timestamps.
$zh.Range("K2").Value = "2016-08-13 12:45:17"
$zh.Range("K3").Value = "2016-08-13 12:45:17"

# The handback is current now, so the "not the latest" error goes away.
$zh.Range("P2").Value = ""
$zh.Range("P3").Value = ""

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")

$de.Columns.Item(3).ColumnWidth = 29.1666666666667
$de.Columns.Item(16).ColumnWidth = 12.8333333333333

$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "Handed back: in sync with en-US"

$de.Range("K2").Value = "2016-08-13 12:45:27"
$de.Range("K3").Value = "2016-08-13 12:45:27"

$de.Range("P2").Value = ""
$de.Range("P3").Value = ""
